$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.848.17'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '1.629.60'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.68'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5105'
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2587'
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06410'
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.40'
$ws.Range('E10').Value = '  -1.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07795'
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.264'
$ws.Range('E12').Value = '  -0.34%  '
$ws.Range('D13').Value = '1.625.90'
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('D14').Value = '1.853.06'
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5592'
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('E16').Value = '  -1.40%  '
$ws.Range('D17').Value = '0.0₅7557'
$ws.Range('E17').Value = '  -2.60%  '
$ws.Range('D18').Value = '25.839.96'
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.57'
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.336'
$ws.Range('E21').Value = '  -2.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.814'
$ws.Range('E22').Value = '  -1.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.009'
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  +0.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.822'
$ws.Range('E25').Value = '  -4.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1294'
$ws.Range('E26').Value = '  +4.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '141.29'
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.755'
$ws.Range('E28').Value = '  -1.45%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.45'
$ws.Range('E29').Value = '  -0.87%  '
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04893'
$ws.Range('E31').Value = '  +0.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.299'
$ws.Range('E32').Value = '  +0.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.192'
$ws.Range('E33').Value = '  -0.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.556'
$ws.Range('E34').Value = '  +1.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.378'
$ws.Range('E35').Value = '  +0.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.8965'
$ws.Range('E36').Value = '  -1.98%  '
$ws.Range('D37').Value = '1.133.71'
$ws.Range('E37').Value = '  -1.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.543'
$ws.Range('E38').Value = '  -1.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5490'
$ws.Range('E39').Value = '  -1.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01558'
$ws.Range('E40').Value = '  -0.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9983'
$ws.Range('E41').Value = '  -0.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.591'
$ws.Range('E42').Value = '  +0.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7958'
$ws.Range('E43').Value = '  -1.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.38'
$ws.Range('E44').Value = '  -2.11%  '
$ws.Range('D45').Value = '1.776.34'
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('D46').Value = '0.0₈112'
$ws.Range('E46').Value = '  -6.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4428'
$ws.Range('E47').Value = '  -2.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.92'
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('E49').Value = '  -2.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.576'
$ws.Range('E50').Value = '  +1.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9999'
$ws.Range('E51').Value = '  -0.59%  '
